$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the population size configuration value
$ws.Range("A2").Value = "Tamanho população: 10000"

# Add new configuration row: number of teachers
$ws.Range("A6").Value = "Número de professores: 18"

# Give column G an explicit width (closest achievable to 43.13 chars)
$ws.Columns.Item(7).ColumnWidth = 42.3

# Move the active selection to A6, matching the edited workbook's view state
[void]$ws.Range("A6").Select()

Write-Output "done"
